$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 39

$ws.Cells.Item($row, 1).Value = "'01/02/2026"
$ws.Cells.Item($row, 1).Style = "Normal"
$ws.Cells.Item($row, 2).Value = 12417.46
$ws.Cells.Item($row, 3).Value = 0.2113242502079662
$ws.Cells.Item($row, 4).Value = 0.7886757497920338
$ws.Cells.Item($row, 5).Value = -133.45
$ws.Cells.Item($row, 6).Value = -22.81
$ws.Cells.Item($row, 7).Value = -20848.41
$ws.Cells.Item($row, 8).Value = -68.04000000000001
$ws.Cells.Item($row, 9).Value = -428.74
$ws.Cells.Item($row, 10).Value = -14.04
